# Add team record (Wins/Losses/Ties) columns to the HOU_1990 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from the last existing header cell (AB1) onto the
# three new header cells so they match the bold/bordered/centered style
# used by the rest of row 1.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every data row (2-46) gets the same team record: 75 wins, 87 losses, 0 ties.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 46) { $lastRow = 46 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 75  # AC
    $ws.Cells.Item($r, 30).Value = 87  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
